$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 39. This shifts the existing rows 39-67 down
# to 40-68 (carrying all of their original values with them), matching the
# weekly roll-forward of this price sheet, and grows the used range to
# A1:R68.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with this week's record.
$ws.Range("A39").Value = 11
$ws.Range("B39").Value = "Vega Monumental Concepción"
$ws.Range("C39").Value = "Bíobío"
$ws.Range("D39").Value = 44582
$ws.Range("E39").Value = 8
$ws.Range("F39").Value = 100112021
$ws.Range("G39").Value = "Ají"
$ws.Range("H39").Value = "Chilena(o)"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 80
$ws.Range("K39").Value = 25000
$ws.Range("L39").Value = 26000
$ws.Range("M39").Value = 25500
$ws.Range("N39").Value = "$/caja 12 kilos"
$ws.Range("O39").Value = "Región Metropolitana"
$ws.Range("P39").Value = 2125
$ws.Range("Q39").Value = 12
$ws.Range("R39").Value = "Hortaliza"
